$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.767.47"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.40%  "
$ws.Range("D3").Value = "'1.634.98"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'215.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.19%  "
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  -0.32%  "
$ws.Range("D9").Value = "'0.0641"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.46%  "
$ws.Range("D10").Value = "'19.85"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.86%  "
$ws.Range("D11").Value = "'0.0778"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.08%  "
$ws.Range("E12").Value = "  -0.93%  "
$ws.Range("D13").Value = "'1.632.25"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.30%  "
$ws.Range("D14").Value = "'1.860.56"
$ws.Range("D14").Style = "Normal"
$ws.Range("E15").Value = "  -1.11%  "
$ws.Range("D17").Value = "'63.14"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("D18").Value = "'25.784.78"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.41%  "
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("E20").Value = "  +2.64%  "
$ws.Range("D21").Value = "'194.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.40%  "
$ws.Range("D22").Value = "'9.96"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.55%  "
$ws.Range("E23").Value = "  +0.47%  "
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("E25").Value = "  -0.61%  "
$ws.Range("D26").Value = "'140.29"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("E27").Value = "  -4.80%  "
$ws.Range("E28").Value = "  +0.19%  "
$ws.Range("E29").Value = "  +0.66%  "
$ws.Range("E30").Value = "  -0.26%  "
$ws.Range("E31").Value = "  +0.71%  "
$ws.Range("E32").Value = "  +1.49%  "
$ws.Range("E33").Value = "  +1.26%  "
$ws.Range("E34").Value = "  +1.48%  "
$ws.Range("E35").Value = "  +0.42%  "
$ws.Range("D36").Value = "'0.898"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.65%  "
$ws.Range("E37").Value = "  -0.26%  "
$ws.Range("D38").Value = "'0.552"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.12%  "
$ws.Range("D39").Value = "'1.107.06"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.78%  "
$ws.Range("E40").Value = "  +0.11%  "
$ws.Range("E41").Value = "  +0.37%  "
$ws.Range("D43").Value = "'0.803"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.14%  "
$ws.Range("D44").Value = "'99.18"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.15%  "
$ws.Range("E45").Value = "  -4.50%  "
$ws.Range("D46").Value = "'55.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.37%  "
$ws.Range("E47").Value = "  +12.21%  "
$ws.Range("E48").Value = "  -2.16%  "
$ws.Range("D49").Value = "'7.70"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.32%  "
$ws.Range("D51").Value = "'0.996"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.65%  "
